$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "weight" column (J) added next to the existing "order" (byte-order) column,
# as part of enhancing the modbus byte-order / register definition.
# Give the new header cell the same look (bold/fill/border) as the other header cells.
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("J1").Value = "weight"

# Fill weight = 1 for each data row (rows 2-16)
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 10).Value = 1
}

# Select the whole new column (matches author's selection after adding the column)
$ws.Range("J1:J1048576").Select()
